# Append: 2025-12-10 01:55 JST
# Update the "取得日時" (acquisition timestamp) column (A) for all data rows
# on the "ランサーズ" sheet from the previous run timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-10 01:24:21"
$newTimestamp = "2025-12-10 01:55:44"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
